# Add a new "Save" column (H) to the s_vals sheet, mirroring the existing
# header style used by the other header cells and filling the new data
# rows with 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the header formatting (bold font, border, centered alignment) from
# an existing header cell (G1) onto the new H1 header cell so it reuses the
# same style record instead of creating a brand-new one.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set the new header text and data values.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
